$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -0.2551564036431285
$ws.Range("C2").Value = 0.7852089611737781
$ws.Range("D2").Value = -0.2475736227641347
$ws.Range("F2").Value = 0.8307123712950741
$ws.Range("G2").Value = 1.100893574250167
$ws.Range("H2").Value = 0.2162804316302916
$ws.Range("I2").Value = -0.4987021728231629
$ws.Range("J2").Value = 1.159308017342042
$ws.Range("K2").Value = 0.5374031615669816
$ws.Range("B3").Value = 0.8076369408067867
$ws.Range("C3").Value = -0.244970083201686
$ws.Range("E3").Value = 0.8247303909907191
$ws.Range("F3").Value = 1.095550324067756
$ws.Range("G3").Value = 0.2059118982599253
$ws.Range("H3").Value = -0.5071794949467867
$ws.Range("I3").Value = 1.151535709035231
$ws.Range("J3").Value = 0.528779188034013
$ws.Range("K3").Value = 0.7685186244804663
$ws.Range("B4").Value = -0.2235455117444659
$ws.Range("D4").Value = 0.7355360465769361
$ws.Range("E4").Value = 1.065823320219178
$ws.Range("F4").Value = 0.2145380182898354
$ws.Range("G4").Value = -0.5331196683722219
$ws.Range("H4").Value = 1.128757146007004
$ws.Range("I4").Value = 0.5151611528853952
$ws.Range("J4").Value = 0.7499978272952905
$ws.Range("C5").Value = 0.7063315727212027
$ws.Range("D5").Value = 1.062802720414257
$ws.Range("E5").Value = 0.195676734678355
$ws.Range("F5").Value = -0.5534507994374261
$ws.Range("G5").Value = 1.114226990572756
$ws.Range("H5").Value = 0.4985201073931272
$ws.Range("I5").Value = 0.7324139818665403
$ws.Range("K5").Value = -0.16111618316075
$ws.Range("B6").Value = 1.046529313339113
$ws.Range("C6").Value = 1.138173782198194
$ws.Range("D6").Value = 0.004583448365355902
$ws.Range("E6").Value = -0.5287568868506929
$ws.Range("F6").Value = 1.123606394129052
$ws.Range("G6").Value = 0.4437093087817362
$ws.Range("H6").Value = 0.7135198909778705
$ws.Range("J6").Value = -0.1898278937901207
$ws.Range("K6").Value = 0.4254186206066807
$ws.Range("B7").Value = 1.588868654188444
$ws.Range("C7").Value = 0.05169415057771429
$ws.Range("D7").Value = -0.7690401176415012
$ws.Range("E7").Value = 1.159551475194915
$ws.Range("F7").Value = 0.4416585513421433
$ws.Range("G7").Value = 0.6369516515418403
$ws.Range("I7").Value = -0.2161169439878461
$ws.Range("J7").Value = 0.382021329893348
$ws.Range("K7").Value = 0.2305062539156956
$ws.Range("B8").Value = 0.3640197375012527
$ws.Range("C8").Value = -0.6357513779059769
$ws.Range("D8").Value = 0.9804730533787229
$ws.Range("E8").Value = 0.4698527901127105
$ws.Range("F8").Value = 0.6729769627342338
$ws.Range("H8").Value = -0.2216592830752073
$ws.Range("I8").Value = 0.3861966588320966
$ws.Range("J8").Value = 0.2193215401759246
$ws.Range("B9").Value = -0.4001608867981357
$ws.Range("C9").Value = 1.065151794253032
$ws.Range("D9").Value = 0.3242524234789934
$ws.Range("E9").Value = 0.6829307055934764
$ws.Range("G9").Value = -0.2605076180326304
$ws.Range("H9").Value = 0.3703508498049107
$ws.Range("I9").Value = 0.2109873117084238
$ws.Range("B10").Value = 1.376175980734415
$ws.Range("C10").Value = 0.4413411978901886
$ws.Range("D10").Value = 0.5204665641439096
$ws.Range("F10").Value = -0.2247647455777268
$ws.Range("G10").Value = 0.3390232034212837
$ws.Range("H10").Value = 0.2077622620068982
$ws.Range("B11").Value = 0.688166375294303
$ws.Range("C11").Value = 0.5381224316817991
$ws.Range("E11").Value = -0.1925473573292467
$ws.Range("F11").Value = 0.351411329249001
$ws.Range("G11").Value = 0.1868984584576193
$ws.Range("B12").Value = 0.777798791098011
$ws.Range("D12").Value = -0.3090319356872608
$ws.Range("E12").Value = 0.3665166009008535
$ws.Range("F12").Value = 0.2101374940836094
$ws.Range("C13").Value = -0.2954946634393446
$ws.Range("D13").Value = 0.3035436819108904
$ws.Range("E13").Value = 0.2201756597651073
$ws.Range("B14").Value = -0.04164562157393659
$ws.Range("C14").Value = 0.4028617320929269
$ws.Range("D14").Value = 0.1085991175498651
$ws.Range("B15").Value = 0.4469214233323758
$ws.Range("C15").Value = 0.130019622424466
$ws.Range("B16").Value = 0.3662627537369125
